$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 640, shifting existing rows 640:681 down to 641:682
$ws.Rows.Item(640).Insert()

# Populate the newly inserted row 640 with the new data point.
# Force column A to text format so the date-like string isn't auto-converted
# into a date serial number by Excel's input parsing, then restore the
# default (General/Normal) style so no stray number-format style is left
# on the cell (matching the plain, unstyled data cells elsewhere).
$ws.Cells.Item(640, 1).NumberFormat = "@"
$ws.Cells.Item(640, 1).Value = "2026/01/17"
$ws.Cells.Item(640, 1).NumberFormat = "General"
$ws.Cells.Item(640, 1).Style = "Normal"
$ws.Cells.Item(640, 2).Value = "土"
$ws.Cells.Item(640, 3).Value = 5
$ws.Cells.Item(640, 4).Value = 35
